$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "akjdhkjfhakhdkjfhaf" text that lived in B7
$ws.Range("B7").ClearContents()

# Add the two new entries (written in this order so the shared-string
# table indices line up with the canonical OOXML: "lljlkfljkajlfjlajf"
# lands before "2+3=5")
$ws.Range("G17").Value = "lljlkfljkajlfjlajf"
$ws.Range("D16").Value = "2+3=5"

# Update the active selection to match the author's final cursor position
$ws.Range("E20").Select()
